# Fixed a bug in WinWeight.merge
# The merge routine produced rows 2-23 in the wrong order; re-sequence
# them back into the correct merge order while keeping the header (row 1)
# and the trailing total rows untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the current (pre-fix) contents of the data block (rows 2-23,
# columns A-F) before overwriting anything.
$original = @{}
for ($r = 2; $r -le 23; $r++) {
    $rowVals = @{}
    for ($c = 1; $c -le 6; $c++) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $original[$r] = $rowVals
}

# Mapping of new row number -> source row number (where its correct data
# currently lives), reflecting the corrected WinWeight.merge order.
$map = @{
    2  = 5
    3  = 13
    4  = 14
    5  = 7
    6  = 15
    7  = 9
    8  = 6
    9  = 2
    10 = 4
    11 = 11
    12 = 8
    13 = 10
    14 = 3
    15 = 12
    16 = 20
    17 = 21
    18 = 17
    19 = 16
    20 = 19
    21 = 18
    22 = 23
    23 = 22
}

foreach ($destRow in $map.Keys) {
    $srcRow = $map[$destRow]
    $srcVals = $original[$srcRow]
    for ($c = 1; $c -le 6; $c++) {
        $ws.Cells.Item($destRow, $c).Value = $srcVals[$c]
    }
}
